$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bug report")
$lo = $ws.ListObjects.Item(1)

# Extend the "Table3" ListObject by four rows (A1:M16 -> A1:M20),
# matching the new Bug Report entries (B016-B019).
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null

# Fill content column-major (B,C,D,E across the four new rows, then
# A last) so new shared-string entries land in the same order as the
# source edit: Search.. / Search / No category show / All categories
# listed / Add with uppercase slug / ... / B016.. B019.
$ws.Range("B17").Value = "Search with space keyword"
$ws.Range("C17").Value = "Search"
$ws.Range("D17").Value = "No category show"
$ws.Range("E17").Value = "All categories listed"

$ws.Range("B18").Value = "Add with uppercase slug"
$ws.Range("C18").Value = "Submit form"
$ws.Range("D18").Value = "Error: Slug cannot contain spaces"
$ws.Range("E18").Value = "Error: Slug must be lowercase"

$ws.Range("B19").Value = "Add category with existing name, different slug"
$ws.Range("C19").Value = "Submit form"
$ws.Range("D19").Value = "Error message displayed"
$ws.Range("E19").Value = "Error: Name must be unique"

$ws.Range("B20").Value = "Add with name containing newline"
$ws.Range("C20").Value = "Submit"
$ws.Range("D20").Value = "Category added successfully"
$ws.Range("E20").Value = "Likely error due to newline"

$ws.Range("A17").Value = "B016"
$ws.Range("A18").Value = "B017"
$ws.Range("A19").Value = "B018"
$ws.Range("A20").Value = "B019"

# Remaining columns: Priority (F), Function ID (G), Severity (H),
# Affected Feature / Version (I), Reported By (J), Status (L) -
# all reuse existing shared strings.
$ws.Range("F17").Value = "Medium"
$ws.Range("F18").Value = "High"
$ws.Range("F19").Value = "High"
$ws.Range("F20").Value = "High"

$ws.Range("G17").Value = "UC02"
$ws.Range("G18").Value = "UC02"
$ws.Range("G19").Value = "UC02"
$ws.Range("G20").Value = "UC02"

$ws.Range("H17").Value = "Normal"
$ws.Range("H18").Value = "Major"
$ws.Range("H19").Value = "Major"
$ws.Range("H20").Value = "Major"

$ws.Range("I17").Value = "Category Management / sprint5-with-bugs"
$ws.Range("I18").Value = "Category Management / sprint5-with-bugs"
$ws.Range("I19").Value = "Category Management / sprint5-with-bugs"
$ws.Range("I20").Value = "Category Management / sprint5-with-bugs"

$ws.Range("J17").Value = "Tran Thi Cat Tuong"
$ws.Range("J18").Value = "Tran Thi Cat Tuong"
$ws.Range("J19").Value = "Tran Thi Cat Tuong"
$ws.Range("J20").Value = "Tran Thi Cat Tuong"

$ws.Range("L17").Value = "Open"
$ws.Range("L18").Value = "Open"
$ws.Range("L19").Value = "Open"
$ws.Range("L20").Value = "Open"

# Date Reported (K): copy the existing date-formatted cell so the new
# cells keep the m/d/yyyy number format (style) instead of minting a
# new one, then write the same serial date (2025-06-11) used by every
# other bug row.
$ws.Range("K16").Copy() | Out-Null
$ws.Range("K17:K20").PasteSpecial(-4122) | Out-Null
$ws.Range("K17").Value = 45819
$ws.Range("K18").Value = 45819
$ws.Range("K19").Value = 45819
$ws.Range("K20").Value = 45819
$excel.CutCopyMode = $false

# Row height to match the other wrapped-text bug rows (42.75pt).
$ws.Rows.Item(17).RowHeight = 42.75
$ws.Rows.Item(18).RowHeight = 42.75
$ws.Rows.Item(19).RowHeight = 42.75
$ws.Rows.Item(20).RowHeight = 42.75

# Column J (Function ID) was widened slightly in the source edit.
$ws.Columns.Item(10).ColumnWidth = 15.333333333333334

# Selection / scroll position after the edit.
$ws.Range("E21").Select()
